# Insert a new weekly price record for "Macroferia Regional de Talca - Repollo"
# (Fruta / hortaliza, semanal) as a new row 185, pushing the existing rows
# 185-196 down to 186-197.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 185:196 down by inserting a new blank row at 185.
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A185").Value = 5
$ws.Range("B185").Value = "Macroferia Regional de Talca"
$ws.Range("C185").Value = "Maule"
$ws.Range("D185").Value = 44516
$ws.Range("E185").Value = 7
$ws.Range("F185").Value = 100112006
$ws.Range("G185").Value = "Repollo"
$ws.Range("H185").Value = "Crespo record"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 3000
$ws.Range("K185").Value = 900
$ws.Range("L185").Value = 900
$ws.Range("M185").Value = 900
$ws.Range("N185").Value = "`$/unidad"
$ws.Range("O185").Value = "Región del Maule"
$ws.Range("P185").Value = 900
$ws.Range("Q185").Value = 1
$ws.Range("R185").Value = "Hortaliza"
